# "add a minigame form" - update the mini-game names (column B / "Name") to
# their localized display text and keep the icon path column (D) pointing at
# the same GameButtonN icons (their shared-string ids simply shift because
# the old, now-unused English game names are dropped from the string table).
#
# Values are written in the same order the original author's commit shows
# the new shared strings appearing in (so the shared string table comes out
# in the same order): B5, B4, B7, B8, B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "魔幻符号"
$ws.Range("B4").Value = "烹饪"
$ws.Range("B7").Value = "三体"
$ws.Range("B8").Value = "SEVEN"
$ws.Range("B6").Value = "智械危机"

# Narrow the "Name" column now that the localized labels are shorter.
$ws.Columns.Item(2).ColumnWidth = 11.375

# Leave the selection on B6, matching the saved workbook view.
$ws.Range("B6").Select()
